# Update the "F" column (sales/attendance count) for a handful of 杭州 events
# on both the "展览" sheet and the "全部类型" aggregate sheet, per the site's
# latest generated snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (row, newValue) updates for column F.
$updates = @{
    "展览"   = @(
        @{ Row = 6;  Value = 5498 }
        @{ Row = 8;  Value = 690 }
        @{ Row = 9;  Value = 948 }
        @{ Row = 13; Value = 587 }
        @{ Row = 17; Value = 1841 }
        @{ Row = 18; Value = 1472 }
        @{ Row = 19; Value = 912 }
        @{ Row = 22; Value = 332 }
        @{ Row = 28; Value = 2893 }
        @{ Row = 40; Value = 719 }
        @{ Row = 41; Value = 88 }
        @{ Row = 44; Value = 68 }
    )
    "全部类型" = @(
        @{ Row = 7;  Value = 5498 }
        @{ Row = 9;  Value = 690 }
        @{ Row = 12; Value = 948 }
        @{ Row = 18; Value = 587 }
        @{ Row = 23; Value = 1841 }
        @{ Row = 24; Value = 1472 }
        @{ Row = 25; Value = 912 }
        @{ Row = 27; Value = 332 }
        @{ Row = 32; Value = 2893 }
        @{ Row = 43; Value = 719 }
        @{ Row = 44; Value = 88 }
        @{ Row = 46; Value = 68 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $ws.Cells.Item($entry.Row, 6).Value = $entry.Value
    }
}
